$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 9-13 (old trailing data no longer present in the scrape)
$ws.Range("A9:H13").EntireRow.Delete()

# Update column widths for columns C, D, F, G, H.
# The host converts the "characters" ColumnWidth into a pixel-rounded raw
# OOXML width, which is offset from the simple integer value. Empirically,
# an input of (target - 1) + 1/12 (plus a small safety margin) lands exactly
# on the desired integer "width" attribute once re-serialized.
$ws.Range("C1").ColumnWidth = 30.103333333333333
$ws.Range("D1").ColumnWidth = 83.10333333333333
$ws.Range("F1").ColumnWidth = 14.103333333333333
$ws.Range("G1").ColumnWidth = 14.103333333333333
$ws.Range("H1").ColumnWidth = 14.103333333333333

# Opportunity IDs in column A are stored as text (not numbers) in the source
# data, so force the Text number format before writing, then restore the
# default "Normal" style so no extra style index lingers on the cells.
$idRange = $ws.Range("A2:A8")
$idRange.NumberFormat = "@"

# Update row 2
$ws.Range("A2").Value = "1328276"
$ws.Range("B2").Value = "https://aiesec.org/opportunity/global-talent/1328276"
$ws.Range("C2").Value = "Web Developer"
$ws.Range("D2").Value = "Tunis, Tunisia"
$ws.Range("F2").Value = "9 applicants"
$ws.Range("H2").Value = "La fabrique"

# Update row 3
$ws.Range("A3").Value = "1328274"
$ws.Range("B3").Value = "https://aiesec.org/opportunity/global-talent/1328274"
$ws.Range("C3").Value = "Web Designer"
$ws.Range("D3").Value = "Tunis, Tunisia"
$ws.Range("F3").Value = "2 applicants"
$ws.Range("H3").Value = "La fabrique"

# Update row 4
$ws.Range("A4").Value = "1328273"
$ws.Range("B4").Value = "https://aiesec.org/opportunity/global-talent/1328273"
$ws.Range("C4").Value = "Digital Marketing Specialist"
$ws.Range("D4").Value = "Tunis, Tunisia"
$ws.Range("F4").Value = "9 applicants"
$ws.Range("H4").Value = "La fabrique"

# Update row 5
$ws.Range("A5").Value = "1328272"
$ws.Range("B5").Value = "https://aiesec.org/opportunity/global-talent/1328272"
$ws.Range("C5").Value = "Graphic Designer"
$ws.Range("D5").Value = "Tunis, Tunisia"
$ws.Range("G5").Value = "9 - 12 Weeks"
$ws.Range("H5").Value = "La fabrique"

# Update row 6
$ws.Range("A6").Value = "1328271"
$ws.Range("B6").Value = "https://aiesec.org/opportunity/global-talent/1328271"
$ws.Range("C6").Value = "Video Editor"
$ws.Range("D6").Value = "Tunis, Tunisia"
$ws.Range("F6").Value = "0 applicants"
$ws.Range("G6").Value = "9 - 12 Weeks"
$ws.Range("H6").Value = "La fabrique"

# Update row 7
$ws.Range("A7").Value = "1327889"
$ws.Range("B7").Value = "https://aiesec.org/opportunity/global-talent/1327889"
$ws.Range("C7").Value = "Graphic Designer"
$ws.Range("D7").Value = "Birkat as SAB, Madinet Berkat as Sabee, Birket el Sab, Menofia Governorate, Egypt"
$ws.Range("F7").Value = "0 applicants"
$ws.Range("H7").Value = "Lines"

# Update row 8
$ws.Range("A8").Value = "1327208"
$ws.Range("B8").Value = "https://aiesec.org/opportunity/global-talent/1327208"
$ws.Range("C8").Value = "Club Manager"
$ws.Range("D8").Value = "Sousse, Tunisia"
$ws.Range("F8").Value = "4 applicants"
$ws.Range("H8").Value = "Saladin"

# Restore the default style on column A so no stray style index is left behind
$idRange.Style = "Normal"
